# TC01_Verify_HomePage.xlsx edit
# Inserts a new "WAIT" step as row 3 on the TC01_Verify_HomePage sheet,
# pushing the existing VERIFY_WEBELEMENT_PRESENT steps (KamanLogo,
# Herobanner, SearchBoxHomePage) down one row (old rows 3-5 become rows 4-6).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC01_Verify_HomePage")

# Insert a blank row above the current row 3; this shifts rows 3:5
# (KamanLogo / Herobanner / SearchBoxHomePage) down to rows 4:6 while
# keeping their values and cell formatting intact.
$ws.Rows("3:3").Insert()

# The newly inserted row 3 has no formatting yet - give it the same
# look (borders/font) as the rows around it by copying the format of
# the row right below it (the former row 3, now row 4).
$ws.Range("A4:E4").Copy($ws.Range("A3:E3"))

# Row 3 becomes the new "WAIT" step: only column B holds a keyword,
# the rest of the row is left blank.
$ws.Range("B3").Value = "WAIT"
$ws.Range("C3:E3").ClearContents()

$excel.CutCopyMode = $false
$ws.Range("B3").Select()
